# Auto-generated edit script for resum_diari_meteocat update
# Commit: Update automàtic: dades i banners [2026-02-08 17:20]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Direct text/value updates (dates, measurements with units, temperatures) ---
# These assign safely as literal text without triggering Excel's numeric/date
# auto-conversion, since they are not pure numeric-looking tokens.
$ws.Range("E2").Value = '2026-02-08 17:18:17'
$ws.Range("I2").Value = '2.7 mm'
$ws.Range("E3").Value = '2026-02-08 17:18:19'
$ws.Range("I3").Value = '1.7 mm'
$ws.Range("K3").Value = '5.8 MJ/m2'
$ws.Range("E4").Value = '2026-02-08 17:18:21'
$ws.Range("J4").Value = '1001.4 hPa'
$ws.Range("E5").Value = '2026-02-08 17:18:24'
$ws.Range("G5").Value = '117 cm'
$ws.Range("I5").Value = '4.9 mm'
$ws.Range("E6").Value = '2026-02-08 17:18:27'
$ws.Range("J6").Value = '1001.3 hPa'
$ws.Range("E7").Value = '2026-02-08 17:18:29'
$ws.Range("J7").Value = '1001.7 hPa'
$ws.Range("K7").Value = '11.3 MJ/m2'
$ws.Range("E8").Value = '2026-02-08 17:18:31'
$ws.Range("J8").Value = '1001.6 hPa'
$ws.Range("K8").Value = '10.9 MJ/m2'
$ws.Range("O8").Value = '9.4 °C'
$ws.Range("E9").Value = '2026-02-08 17:18:34'
$ws.Range("E10").Value = '2026-02-08 17:18:36'
$ws.Range("K10").Value = '11.3 MJ/m2'
$ws.Range("O10").Value = '9.3 °C'
$ws.Range("E11").Value = '2026-02-08 17:18:39'
$ws.Range("O11").Value = '4.6 °C'
$ws.Range("E12").Value = '2026-02-08 17:18:41'
$ws.Range("E13").Value = '2026-02-08 17:18:43'
$ws.Range("O13").Value = '3.3 °C'
$ws.Range("E14").Value = '2026-02-08 17:18:46'
$ws.Range("K14").Value = '11.9 MJ/m2'
$ws.Range("O14").Value = '11.0 °C'
$ws.Range("E15").Value = '2026-02-08 17:18:48'
$ws.Range("E16").Value = '2026-02-08 17:18:51'
$ws.Range("I16").Value = '2.0 mm'
$ws.Range("K16").Value = '7.6 MJ/m2'
$ws.Range("E17").Value = '2026-02-08 17:18:53'
$ws.Range("E18").Value = '2026-02-08 17:18:55'
$ws.Range("J18").Value = '1001.7 hPa'
$ws.Range("E19").Value = '2026-02-08 17:18:58'
$ws.Range("E20").Value = '2026-02-08 17:19:00'
$ws.Range("E21").Value = '2026-02-08 17:19:02'
$ws.Range("K21").Value = '11.6 MJ/m2'
$ws.Range("O21").Value = '5.2 °C'
$ws.Range("E22").Value = '2026-02-08 17:19:05'
$ws.Range("K22").Value = '7.5 MJ/m2'
$ws.Range("E23").Value = '2026-02-08 17:19:07'
$ws.Range("I23").Value = '3.7 mm'
$ws.Range("E24").Value = '2026-02-08 17:19:10'
$ws.Range("J24").Value = '1003.1 hPa'
$ws.Range("K24").Value = '11.4 MJ/m2'
$ws.Range("O24").Value = '8.3 °C'
$ws.Range("E25").Value = '2026-02-08 17:19:12'
$ws.Range("E26").Value = '2026-02-08 17:19:14'
$ws.Range("J26").Value = '1000.6 hPa'
$ws.Range("E27").Value = '2026-02-08 17:19:17'
$ws.Range("O27").Value = '-3.0 °C'
$ws.Range("E28").Value = '2026-02-08 17:19:19'
$ws.Range("E29").Value = '2026-02-08 17:19:22'
$ws.Range("K29").Value = '11.8 MJ/m2'
$ws.Range("E30").Value = '2026-02-08 17:19:24'
$ws.Range("J30").Value = '1001.7 hPa'
$ws.Range("K30").Value = '10.4 MJ/m2'
$ws.Range("E31").Value = '2026-02-08 17:19:26'
$ws.Range("J31").Value = '1000.7 hPa'
$ws.Range("N31").Value = '8.9 °C 16:58 TU'
$ws.Range("E32").Value = '2026-02-08 17:19:28'
$ws.Range("O32").Value = '3.8 °C'
$ws.Range("E33").Value = '2026-02-08 17:19:31'
$ws.Range("J33").Value = '1002.9 hPa'
$ws.Range("E34").Value = '2026-02-08 17:19:33'
$ws.Range("K34").Value = '12.6 MJ/m2'
$ws.Range("E35").Value = '2026-02-08 17:19:36'
$ws.Range("J35").Value = '1003.7 hPa'
$ws.Range("E36").Value = '2026-02-08 17:19:38'
$ws.Range("J36").Value = '1001.7 hPa'
$ws.Range("K36").Value = '11.5 MJ/m2'
$ws.Range("E37").Value = '2026-02-08 17:19:41'
$ws.Range("O37").Value = '5.4 °C'
$ws.Range("E38").Value = '2026-02-08 17:19:43'
$ws.Range("E39").Value = '2026-02-08 17:19:46'
$ws.Range("K39").Value = '12.8 MJ/m2'
$ws.Range("E40").Value = '2026-02-08 17:19:48'
$ws.Range("J40").Value = '1003.2 hPa'
$ws.Range("O40").Value = '5.8 °C'
$ws.Range("E41").Value = '2026-02-08 17:19:51'
$ws.Range("J41").Value = '1001.9 hPa'
$ws.Range("K41").Value = '12.6 MJ/m2'
$ws.Range("O41").Value = '11.8 °C'
$ws.Range("E42").Value = '2026-02-08 17:19:53'
$ws.Range("E43").Value = '2026-02-08 17:19:55'
$ws.Range("O43").Value = '6.6 °C'
$ws.Range("E44").Value = '2026-02-08 17:19:58'
$ws.Range("I44").Value = '1.6 mm'
$ws.Range("K44").Value = '6.1 MJ/m2'
$ws.Range("E45").Value = '2026-02-08 17:20:00'
$ws.Range("I45").Value = '1.6 mm'
$ws.Range("J45").Value = '1003.8 hPa'
$ws.Range("E46").Value = '2026-02-08 17:20:03'
$ws.Range("J46").Value = '1003.5 hPa'
$ws.Range("K46").Value = '8.0 MJ/m2'
$ws.Range("O46").Value = '9.2 °C'

# --- Percentage-looking text updates (e.g. "78%") ---
# Assigning a bare "NN%" string straight to a General-formatted cell makes Excel
# "smart"-convert it into a numeric percentage (changing both the stored type and
# the cell style). To keep these as literal text in their original style, stage the
# text in a scratch cell explicitly formatted as Text ("@"), copy it, and paste
# (Values only) into the destination - this preserves the destination's existing
# style/number format while still storing the literal string.
$scratch = $ws.Range("AA1")
$scratch.NumberFormat = "@"

$scratch.Value = '78%'
$scratch.Copy()
$ws.Range("H13").PasteSpecial(-4163)

$scratch.Value = '77%'
$scratch.Copy()
$ws.Range("H14").PasteSpecial(-4163)

$scratch.Value = '69%'
$scratch.Copy()
$ws.Range("H18").PasteSpecial(-4163)

$scratch.Value = '82%'
$scratch.Copy()
$ws.Range("H21").PasteSpecial(-4163)

$scratch.Value = '83%'
$scratch.Copy()
$ws.Range("H24").PasteSpecial(-4163)

$scratch.Value = '68%'
$scratch.Copy()
$ws.Range("H30").PasteSpecial(-4163)

$scratch.Value = '84%'
$scratch.Copy()
$ws.Range("H40").PasteSpecial(-4163)

$scratch.Value = '87%'
$scratch.Copy()
$ws.Range("H43").PasteSpecial(-4163)

$scratch.Value = '76%'
$scratch.Copy()
$ws.Range("H46").PasteSpecial(-4163)

# Clean up the scratch cell so it leaves no trace in the saved workbook
$scratch.ClearContents()
$scratch.ClearFormats()
$excel.CutCopyMode = $false

